$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current last data row (row 69),
# shifting the existing row 69 down to row 71.
$ws.Rows("69:70").Insert()

# Row 69 - new weekly entry, Primera quality
$ws.Range("A69").Value = 7
$ws.Range("B69").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C69").Value = "Ñuble"
$ws.Range("D69").Value = 44783
$ws.Range("E69").Value = 16
$ws.Range("F69").Value = 100112040
$ws.Range("G69").Value = "Cilantro"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 200
$ws.Range("K69").Value = 700
$ws.Range("L69").Value = 800
$ws.Range("M69").Value = 750
$ws.Range("N69").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O69").Value = "Provincia de Diguillín"
$ws.Range("P69").Value = 750
$ws.Range("Q69").Value = 1
$ws.Range("R69").Value = "Hortaliza"

# Row 70 - new weekly entry, Segunda quality
$ws.Range("A70").Value = 7
$ws.Range("B70").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C70").Value = "Ñuble"
$ws.Range("D70").Value = 44783
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 100112040
$ws.Range("G70").Value = "Cilantro"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Segunda"
$ws.Range("J70").Value = 200
$ws.Range("K70").Value = 600
$ws.Range("L70").Value = 600
$ws.Range("M70").Value = 600
$ws.Range("N70").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O70").Value = "Provincia de Diguillín"
$ws.Range("P70").Value = 600
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"

# Ensure the date cells use the same date/time number format as the rest of column D
$ws.Range("D69:D70").NumberFormat = $ws.Range("D71").NumberFormat
